# Auto-generated edit script: update cryptos price/volume columns
# D = Price (forced as text via leading apostrophe to avoid numeric coercion)
# E = Volume(1h) percentage strings (already non-numeric text)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '23.404.82'
$ws.Range("E2").Value = '  -1.29%  '
$ws.Range("D3").Value = "'" + '1.636.82'
$ws.Range("E3").Value = '  -0.96%  '
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("E5").Value = '  -0.13%  '
$ws.Range("D6").Value = "'" + '299.99'
$ws.Range("E6").Value = '  -1.06%  '
$ws.Range("D7").Value = "'" + '0.3780'
$ws.Range("E7").Value = '  -0.54%  '
$ws.Range("D8").Value = "'" + '50.29'
$ws.Range("E8").Value = '  -1.34%  '
$ws.Range("D9").Value = "'" + '0.3526'
$ws.Range("E9").Value = '  -2.58%  '
$ws.Range("D10").Value = "'" + '0.08058'
$ws.Range("E10").Value = '  -1.77%  '
$ws.Range("D11").Value = "'" + '1.210'
$ws.Range("E11").Value = '  -3.01%  '
$ws.Range("D12").Value = "'" + '1.001'
$ws.Range("E12").Value = '  -0.20%  '
$ws.Range("D13").Value = "'" + '21.98'
$ws.Range("E13").Value = '  -3.02%  '
$ws.Range("D14").Value = "'" + '6.339'
$ws.Range("E14").Value = '  -2.76%  '
$ws.Range("D15").Value = "'" + '7.259'
$ws.Range("E15").Value = '  -2.53%  '
$ws.Range("E16").Value = '  -2.51%  '
$ws.Range("D17").Value = "'" + '1.640.62'
$ws.Range("E17").Value = '  -0.97%  '
$ws.Range("D18").Value = "'" + '95.99'
$ws.Range("E18").Value = '  -1.30%  '
$ws.Range("D19").Value = "'" + '0.06949'
$ws.Range("E19").Value = '  -0.98%  '
$ws.Range("D20").Value = "'" + '6.687'
$ws.Range("E20").Value = '  -1.45%  '
$ws.Range("D21").Value = "'" + '17.34'
$ws.Range("D22").Value = "'" + '1.001'
$ws.Range("E22").Value = '  -0.15%  '
$ws.Range("E23").Value = '  -4.20%  '
$ws.Range("D24").Value = "'" + '23.430.25'
$ws.Range("E24").Value = '  -1.23%  '
$ws.Range("D25").Value = "'" + '2.468'
$ws.Range("E25").Value = '  -2.48%  '
$ws.Range("D26").Value = "'" + '2.897'
$ws.Range("E26").Value = '  -5.09%  '
$ws.Range("D27").Value = "'" + '20.82'
$ws.Range("E27").Value = '  -2.12%  '
$ws.Range("D28").Value = "'" + '151.93'
$ws.Range("E28").Value = '  +0.28%  '
$ws.Range("D29").Value = "'" + '5.188'
$ws.Range("E29").Value = '  -0.68%  '
$ws.Range("D30").Value = "'" + '132.72'
$ws.Range("E30").Value = '  -1.23%  '
$ws.Range("D31").Value = "'" + '1.817.82'
$ws.Range("E31").Value = '  -1.16%  '
$ws.Range("D32").Value = "'" + '6.841'
$ws.Range("E32").Value = '  -1.22%  '
$ws.Range("D33").Value = "'" + '2.129'
$ws.Range("E33").Value = '  -4.37%  '
$ws.Range("E34").Value = '  -3.24%  '
$ws.Range("D35").Value = "'" + '0.9769'
$ws.Range("E35").Value = '  -9.07%  '
$ws.Range("D36").Value = "'" + '0.02704'
$ws.Range("E36").Value = '  -3.67%  '
$ws.Range("D37").Value = "'" + '0.08720'
$ws.Range("D39").Value = "'" + '5.883'
$ws.Range("E39").Value = '  -3.57%  '
$ws.Range("D40").Value = "'" + '12.98'
$ws.Range("E40").Value = '  -0.26%  '
$ws.Range("D41").Value = "'" + '0.06799'
$ws.Range("E41").Value = '  -3.91%  '
$ws.Range("D42").Value = "'" + '0.6843'
$ws.Range("E42").Value = '  -2.64%  '
$ws.Range("E43").Value = '  -2.64%  '
$ws.Range("D44").Value = "'" + '15.62'
$ws.Range("E44").Value = '  -1.96%  '
$ws.Range("E45").Value = '  -0.01%  '
$ws.Range("D46").Value = "'" + '0.6323'
$ws.Range("E46").Value = '  -2.77%  '
$ws.Range("E47").Value = '  -3.44%  '
$ws.Range("D48").Value = "'" + '3.903'
$ws.Range("E48").Value = '  -1.43%  '
$ws.Range("D49").Value = "'" + '0.07695'
$ws.Range("E49").Value = '  -3.28%  '
$ws.Range("D50").Value = "'" + '126.90'
$ws.Range("D51").Value = "'" + '1.218'
$ws.Range("E51").Value = '  +2.06%  '
